$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.555.77'
$ws.Range('E2').Value = '  -2.98%  '
$ws.Range('D3').Value = '2.363.81'
$ws.Range('E3').Value = '  -4.63%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '309.62'
$ws.Range('E5').Value = '  -2.93%  '
$ws.Range('D6').Value = '86.85'
$ws.Range('E6').Value = '  -7.05%  '
$ws.Range('D7').Value = '0.527'
$ws.Range('E7').Value = '  -4.84%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '0.493'
$ws.Range('E9').Value = '  -5.15%  '
$ws.Range('D10').Value = '0.0836'
$ws.Range('E10').Value = '  -5.36%  '
$ws.Range('D11').Value = '30.58'
$ws.Range('E11').Value = '  -8.15%  '
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '2.743.94'
$ws.Range('E13').Value = '  -4.05%  '
$ws.Range('D14').Value = '6.52'
$ws.Range('E14').Value = '  -6.01%  '
$ws.Range('D15').Value = '14.91'
$ws.Range('E15').Value = '  -5.17%  '
$ws.Range('D16').Value = '2.376.92'
$ws.Range('E16').Value = '  -3.61%  '
$ws.Range('D17').Value = '0.756'
$ws.Range('E17').Value = '  -5.96%  '
$ws.Range('D18').Value = '40.542.24'
$ws.Range('E18').Value = '  -2.89%  '
$ws.Range('D19').Value = '0.0₃0905'
$ws.Range('E19').Value = '  -5.23%  '
$ws.Range('D20').Value = '6.11'
$ws.Range('E20').Value = '  -5.89%  '
$ws.Range('D21').Value = '68.52'
$ws.Range('E21').Value = '  -3.89%  '
$ws.Range('D22').Value = '10.75'
$ws.Range('E22').Value = '  -5.44%  '
$ws.Range('D23').Value = '232.03'
$ws.Range('E23').Value = '  -4.08%  '
$ws.Range('D24').Value = '2.62'
$ws.Range('E24').Value = '  -5.13%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E26').Value = '  -7.77%  '
$ws.Range('D27').Value = '23.63'
$ws.Range('E27').Value = '  -6.67%  '
$ws.Range('D28').Value = '2.20'
$ws.Range('E28').Value = '  -2.86%  '
$ws.Range('D29').Value = '9.30'
$ws.Range('E29').Value = '  -4.80%  '
$ws.Range('D30').Value = '33.67'
$ws.Range('E30').Value = '  -9.29%  '
$ws.Range('D31').Value = '152.87'
$ws.Range('E31').Value = '  -3.81%  '
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = '5.19'
$ws.Range('E33').Value = '  -6.29%  '
$ws.Range('D34').Value = '0.0726'
$ws.Range('E34').Value = '  -5.40%  '
$ws.Range('E35').Value = '  -4.74%  '
$ws.Range('D36').Value = '0.113'
$ws.Range('E36').Value = '  -2.74%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '0.0996'
$ws.Range('E37').Value = '  -4.72%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').Value = '15.78'
$ws.Range('E38').Value = '  -9.68%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '2.75'
$ws.Range('E39').Value = '  -6.28%  '
$ws.Range('D40').Value = '1.70'
$ws.Range('E40').Value = '  -9.66%  '
$ws.Range('D41').Value = '3.85'
$ws.Range('E41').Value = '  -4.64%  '
$ws.Range('E42').Value = '  -4.73%  '
$ws.Range('D43').Value = '1.951.33'
$ws.Range('E43').Value = '  -2.82%  '
$ws.Range('D44').Value = '0.0270'
$ws.Range('E44').Value = '  -5.45%  '
$ws.Range('D45').Value = '17.49'
$ws.Range('E45').Value = '  -9.45%  '
$ws.Range('D46').Value = '9.54'
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('D47').Value = '2.69'
$ws.Range('E47').Value = '  -9.72%  '
$ws.Range('D48').Value = '2.608.66'
$ws.Range('E48').Value = '  -3.97%  '
$ws.Range('D49').Value = '92.88'
$ws.Range('E49').Value = '  -5.58%  '
$ws.Range('D50').Value = '72.10'
$ws.Range('E50').Value = '  -6.26%  '
$ws.Range('D51').Value = '50.38'
$ws.Range('E51').Value = '  -4.13%  '
